$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.072547900429899
$ws.Cells.Item(2, 4).Value = 1.073573566595749
$ws.Cells.Item(2, 5).Value = 1.076313147356533
$ws.Cells.Item(2, 6).Value = 1.085891919818932
$ws.Cells.Item(2, 9).Value = 1.055202800220539
$ws.Cells.Item(2, 10).Value = 1.077466893430117
$ws.Cells.Item(2, 11).Value = 1.076265068450713
$ws.Cells.Item(2, 12).Value = 1.078997408825591
$ws.Cells.Item(2, 13).Value = 1.088551175733921
$ws.Cells.Item(2, 14).Value = 1.078997020005968
# Row 3
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.073919608752115
$ws.Cells.Item(3, 4).Value = 1.074503739208428
$ws.Cells.Item(3, 5).Value = 1.07750022718533
$ws.Cells.Item(3, 6).Value = 1.087074793978532
$ws.Cells.Item(3, 9).Value = 1.05559230270504
$ws.Cells.Item(3, 10).Value = 1.078494969104629
$ws.Cells.Item(3, 11).Value = 1.077011923151688
$ws.Cells.Item(3, 12).Value = 1.080001065207994
$ws.Cells.Item(3, 13).Value = 1.089552451524885
$ws.Cells.Item(3, 14).Value = 1.080026555665859
# Row 4
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.074806752152627
$ws.Cells.Item(4, 4).Value = 1.075105222314382
$ws.Cells.Item(4, 5).Value = 1.078268172485991
$ws.Cells.Item(4, 6).Value = 1.08784009917592
$ws.Cells.Item(4, 9).Value = 1.055842896935975
$ws.Cells.Item(4, 10).Value = 1.079159254354905
$ws.Cells.Item(4, 11).Value = 1.077494156752617
$ws.Cells.Item(4, 12).Value = 1.080649743858838
$ws.Cells.Item(4, 13).Value = 1.090199667415747
$ws.Cells.Item(4, 14).Value = 1.080691784277386
# Row 5
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.075179604664979
$ws.Cells.Item(5, 4).Value = 1.075357991004961
$ws.Cells.Item(5, 5).Value = 1.078590977126341
$ws.Cells.Item(5, 6).Value = 1.088161813094164
$ws.Cells.Item(5, 9).Value = 1.055947902905047
$ws.Cells.Item(5, 10).Value = 1.079438295212492
$ws.Cells.Item(5, 11).Value = 1.077696642095048
$ws.Cells.Item(5, 12).Value = 1.080922269930159
$ws.Cells.Item(5, 13).Value = 1.090471596887072
$ws.Cells.Item(5, 14).Value = 1.080971221404997
# Row 6
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.075242202332827
$ws.Cells.Item(6, 4).Value = 1.075400426453796
$ws.Cells.Item(6, 5).Value = 1.07864517519157
$ws.Cells.Item(6, 6).Value = 1.088215829138977
$ws.Cells.Item(6, 9).Value = 1.055965513731415
$ws.Cells.Item(6, 10).Value = 1.079485134295742
$ws.Cells.Item(6, 11).Value = 1.077730625905138
$ws.Cells.Item(6, 12).Value = 1.0809680178545
$ws.Cells.Item(6, 13).Value = 1.090517245708292
$ws.Cells.Item(6, 14).Value = 1.081018127005118
# Row 7
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.074811734624664
$ws.Cells.Item(7, 4).Value = 1.075108600193112
$ws.Cells.Item(7, 5).Value = 1.078272485970033
$ws.Cells.Item(7, 6).Value = 1.087844398010666
$ws.Cells.Item(7, 9).Value = 1.055844301380863
$ws.Cells.Item(7, 10).Value = 1.079162983792281
$ws.Cells.Item(7, 11).Value = 1.077496863335424
$ws.Cells.Item(7, 12).Value = 1.08065338606404
$ws.Cells.Item(7, 13).Value = 1.090203301577373
$ws.Cells.Item(7, 14).Value = 1.080695519010991
# Row 8
$ws.Cells.Item(8, 2).Value = 1.019999999999999
$ws.Cells.Item(8, 3).Value = 1.073011569716424
$ws.Cells.Item(8, 4).Value = 1.073888006455632
$ws.Cells.Item(8, 5).Value = 1.076714364086426
$ws.Cells.Item(8, 6).Value = 1.086291698236441
$ws.Cells.Item(8, 9).Value = 1.055334733081303
$ws.Cells.Item(8, 10).Value = 1.077814534248631
$ws.Cells.Item(8, 11).Value = 1.076517685785305
$ws.Cells.Item(8, 12).Value = 1.079336756620368
$ws.Cells.Item(8, 13).Value = 1.088889702774091
$ws.Cells.Item(8, 14).Value = 1.079345154514319
# Row 9
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.069835889049304
$ws.Cells.Item(9, 4).Value = 1.071734043266825
$ws.Cells.Item(9, 5).Value = 1.073967314174465
$ws.Cells.Item(9, 6).Value = 1.083554837173796
$ws.Cells.Item(9, 9).Value = 1.054425741151816
$ws.Cells.Item(9, 10).Value = 1.075431021989251
$ws.Cells.Item(9, 11).Value = 1.074784295101359
$ws.Cells.Item(9, 12).Value = 1.077010816665652
$ws.Cells.Item(9, 13).Value = 1.086569708659989
$ws.Cells.Item(9, 14).Value = 1.076958257394138
# Row 10
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.067716145875075
$ws.Cells.Item(10, 4).Value = 1.070295894530618
$ws.Cells.Item(10, 5).Value = 1.072134834565536
$ws.Cells.Item(10, 6).Value = 1.081729593441593
$ws.Cells.Item(10, 9).Value = 1.053812249059347
$ws.Cells.Item(10, 10).Value = 1.073836904648845
$ws.Cells.Item(10, 11).Value = 1.073623272626794
$ws.Cells.Item(10, 12).Value = 1.075456108270368
$ws.Cells.Item(10, 13).Value = 1.085019386165857
$ws.Cells.Item(10, 14).Value = 1.075361876224261
# Row 11
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.066797600057536
$ws.Cells.Item(11, 4).Value = 1.069672627294882
$ws.Cells.Item(11, 5).Value = 1.071341053735281
$ws.Cells.Item(11, 6).Value = 1.080939054639783
$ws.Cells.Item(11, 9).Value = 1.053544808974473
$ws.Cells.Item(11, 10).Value = 1.073145387808049
$ws.Cells.Item(11, 11).Value = 1.073119230913392
$ws.Cells.Item(11, 12).Value = 1.074781903616518
$ws.Cells.Item(11, 13).Value = 1.084347184009574
$ws.Cells.Item(11, 14).Value = 1.074669377350233
# Row 12
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.066456304099269
$ws.Cells.Item(12, 4).Value = 1.069441035457459
$ws.Cells.Item(12, 5).Value = 1.07104615942919
$ws.Cells.Item(12, 6).Value = 1.080645381231323
$ws.Cells.Item(12, 9).Value = 1.053445199118291
$ws.Cells.Item(12, 10).Value = 1.072888336266364
$ws.Cells.Item(12, 11).Value = 1.072931808679845
$ws.Cells.Item(12, 12).Value = 1.074531320173284
$ws.Cells.Item(12, 13).Value = 1.084097360190085
$ws.Cells.Item(12, 14).Value = 1.074411960765874
# Row 13
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.066529518162168
$ws.Cells.Item(13, 4).Value = 1.069490716458047
$ws.Cells.Item(13, 5).Value = 1.071109417512687
$ws.Cells.Item(13, 6).Value = 1.08070837666906
$ws.Cells.Item(13, 9).Value = 1.053466578039336
$ws.Cells.Item(13, 10).Value = 1.072943483405209
$ws.Cells.Item(13, 11).Value = 1.072972020403894
$ws.Cells.Item(13, 12).Value = 1.074585078176885
$ws.Cells.Item(13, 13).Value = 1.084150954532077
$ws.Cells.Item(13, 14).Value = 1.074467186219982
# Row 14
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.06676939063908
$ws.Cells.Item(14, 4).Value = 1.06965348553641
$ws.Cells.Item(14, 5).Value = 1.071316678691683
$ws.Cells.Item(14, 6).Value = 1.080914780174226
$ws.Cells.Item(14, 9).Value = 1.053536580721737
$ws.Cells.Item(14, 10).Value = 1.073124143769679
$ws.Cells.Item(14, 11).Value = 1.07310374259629
$ws.Cells.Item(14, 12).Value = 1.074761193464542
$ws.Cells.Item(14, 13).Value = 1.084326536324278
$ws.Cells.Item(14, 14).Value = 1.074648103142892
# Row 15
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.06691716965573
$ws.Cells.Item(15, 4).Value = 1.069753761939496
$ws.Cells.Item(15, 5).Value = 1.071444372595162
$ws.Cells.Item(15, 6).Value = 1.081041947853991
$ws.Cells.Item(15, 9).Value = 1.053579675778917
$ws.Cells.Item(15, 10).Value = 1.073235429128711
$ws.Cells.Item(15, 11).Value = 1.073184874606257
$ws.Cells.Item(15, 12).Value = 1.074869683441948
$ws.Cells.Item(15, 13).Value = 1.084434699700371
$ws.Cells.Item(15, 14).Value = 1.074759546539898
# Row 16
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.067777091959183
$ws.Cells.Item(16, 4).Value = 1.070337247224304
$ws.Cells.Item(16, 5).Value = 1.072187508507799
$ws.Cells.Item(16, 6).Value = 1.081782054540908
$ws.Cells.Item(16, 9).Value = 1.053829960279804
$ws.Cells.Item(16, 10).Value = 1.07388277167218
$ws.Cells.Item(16, 11).Value = 1.073656696443433
$ws.Cells.Item(16, 12).Value = 1.075500831592229
$ws.Cells.Item(16, 13).Value = 1.08506397879009
$ws.Cells.Item(16, 14).Value = 1.07540780838403
# Row 17
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.068316312443149
$ws.Cells.Item(17, 4).Value = 1.070703106600565
$ws.Cells.Item(17, 5).Value = 1.072653574292452
$ws.Cells.Item(17, 6).Value = 1.082246249532959
$ws.Cells.Item(17, 9).Value = 1.053986476009724
$ws.Cells.Item(17, 10).Value = 1.074288494781322
$ws.Cells.Item(17, 11).Value = 1.0739523055498
$ws.Cells.Item(17, 12).Value = 1.075896462986653
$ws.Cells.Item(17, 13).Value = 1.085458465643219
$ws.Cells.Item(17, 14).Value = 1.075814107666525
# Row 18
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.068630765033047
$ws.Cells.Item(18, 4).Value = 1.070916454256206
$ws.Cells.Item(18, 5).Value = 1.072925393344529
$ws.Cells.Item(18, 6).Value = 1.082516987864468
$ws.Cells.Item(18, 9).Value = 1.054077595976145
$ws.Cells.Item(18, 10).Value = 1.074525025469778
$ws.Cells.Item(18, 11).Value = 1.074124602932176
$ws.Cells.Item(18, 12).Value = 1.076127131239676
$ws.Cells.Item(18, 13).Value = 1.085688476260093
$ws.Cells.Item(18, 14).Value = 1.076050974255688
# Row 19
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.068737974250115
$ws.Cells.Item(19, 4).Value = 1.070989191533342
$ws.Cells.Item(19, 5).Value = 1.073018071746104
$ws.Cells.Item(19, 6).Value = 1.082609299536611
$ws.Cells.Item(19, 9).Value = 1.054108636219688
$ws.Cells.Item(19, 10).Value = 1.074605655935321
$ws.Cells.Item(19, 11).Value = 1.074183330496574
$ws.Cells.Item(19, 12).Value = 1.076205766808162
$ws.Cells.Item(19, 13).Value = 1.085766889271885
$ws.Cells.Item(19, 14).Value = 1.076131719225741
# Row 20
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.06825846600742
$ws.Cells.Item(20, 4).Value = 1.07066385873569
$ws.Cells.Item(20, 5).Value = 1.072603572890849
$ws.Cells.Item(20, 6).Value = 1.082196447763301
$ws.Cells.Item(20, 9).Value = 1.053969701257528
$ws.Cells.Item(20, 10).Value = 1.074244977043175
$ws.Cells.Item(20, 11).Value = 1.073920602578309
$ws.Cells.Item(20, 12).Value = 1.075854025520272
$ws.Cells.Item(20, 13).Value = 1.085416149943522
$ws.Cells.Item(20, 14).Value = 1.075770528128198
# Row 21
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.066698757167396
$ws.Cells.Item(21, 4).Value = 1.069605556393822
$ws.Cells.Item(21, 5).Value = 1.071255646809999
$ws.Cells.Item(21, 6).Value = 1.080854000393603
$ws.Cells.Item(21, 9).Value = 1.053515974158699
$ws.Cells.Item(21, 10).Value = 1.073070949103408
$ws.Cells.Item(21, 11).Value = 1.073064959192375
$ws.Cells.Item(21, 12).Value = 1.074709336169913
$ws.Cells.Item(21, 13).Value = 1.084274835692193
$ws.Cells.Item(21, 14).Value = 1.074594832934094
# Row 22
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.065717484164404
$ws.Cells.Item(22, 4).Value = 1.068939680803644
$ws.Cells.Item(22, 5).Value = 1.070407868598198
$ws.Cells.Item(22, 6).Value = 1.080009763431424
$ws.Cells.Item(22, 9).Value = 1.053229131144549
$ws.Cells.Item(22, 10).Value = 1.072331682245863
$ws.Cells.Item(22, 11).Value = 1.072525831783419
$ws.Cells.Item(22, 12).Value = 1.07398873325315
$ws.Cells.Item(22, 13).Value = 1.083556446425404
$ws.Cells.Item(22, 14).Value = 1.073854516232815
# Row 23
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.066237736036318
$ws.Cells.Item(23, 4).Value = 1.069292719985474
$ws.Cells.Item(23, 5).Value = 1.070857319730227
$ws.Cells.Item(23, 6).Value = 1.080457328022731
$ws.Cells.Item(23, 9).Value = 1.053381340960467
$ws.Cells.Item(23, 10).Value = 1.072723687732664
$ws.Cells.Item(23, 11).Value = 1.07281174303837
$ws.Cells.Item(23, 12).Value = 1.074370823890319
$ws.Cells.Item(23, 13).Value = 1.083937354790768
$ws.Cells.Item(23, 14).Value = 1.074247078412371
# Row 24
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.06828460452424
$ws.Cells.Item(24, 4).Value = 1.070681593316924
$ws.Cells.Item(24, 5).Value = 1.072626166462229
$ws.Cells.Item(24, 6).Value = 1.08221895109725
$ws.Cells.Item(24, 9).Value = 1.053977281580545
$ws.Cells.Item(24, 10).Value = 1.074264641208806
$ws.Cells.Item(24, 11).Value = 1.073934928177073
$ws.Cells.Item(24, 12).Value = 1.075873201483804
$ws.Cells.Item(24, 13).Value = 1.085435270855204
$ws.Cells.Item(24, 14).Value = 1.075790220219199
# Row 25
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.070657324060693
$ws.Cells.Item(25, 4).Value = 1.072291271914288
$ws.Cells.Item(25, 5).Value = 1.074677677947614
$ws.Cells.Item(25, 6).Value = 1.084262489978139
$ws.Cells.Item(25, 9).Value = 1.054662054692564
$ws.Cells.Item(25, 10).Value = 1.076048106336629
$ws.Cells.Item(25, 11).Value = 1.075233369065401
$ws.Cells.Item(25, 12).Value = 1.077612837429947
$ws.Cells.Item(25, 13).Value = 1.0871701189938117
$ws.Cells.Item(25, 14).Value = 1.07757621807207

Write-Host "Updated vm_pu values for 380 kV case"